# Refactor email code into class
# (Data update: roll rates forward from 2021-05-23 to 2021-05-24 and
#  recompute converted amounts using the new USD rate 1.2035.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: update the date for all data rows (1-12).
# Leading apostrophe forces these to stay plain text instead of being
# auto-converted into an Excel date serial number.
$ws.Range("A1:A12").Value = "'2021-05-24"

# Column G: update the USD conversion rate for rows 1-9.
$ws.Range("G1:G9").Value = "'1.2035"

# Column H: recomputed converted amounts (Column B * new rate), stored
# as text to match the original inline-string formatting.
$ws.Range("H1").Value = "'858.0955"
$ws.Range("H2").Value = "'4427.6765000000005"
$ws.Range("H3").Value = "'28.9039781"
$ws.Range("H4").Value = "'23.588600000000003"
$ws.Range("H5").Value = "'8063.45"
$ws.Range("H6").Value = "'4.440915"
$ws.Range("H7").Value = "'68.21438"
$ws.Range("H8").Value = "'67.85333"
